$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (VieonDpoint / Dpoint@2021 with hyperlink) into rows 3 and 4,
# copying values + number formatting/style from row 2 so the new cells line up
# with the existing "Hyperlink" cell style (s="1").
$ws.Range("A2:B2").Copy($ws.Range("A3:B3"))
$ws.Range("A2:B2").Copy($ws.Range("A4:B4"))

# Re-create the hyperlink on the two new cells (same mailto target as B2).
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Dpoint@2021") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Dpoint@2021") | Out-Null

# Adding a hyperlink re-stamps a freshly derived style on the target cell;
# restore the original "Hyperlink" cell-style formatting (shared with B2)
# on the new cells so they match row 2 exactly.
$ws.Range("B2").Copy()
$ws.Range("B3:B4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection to G4, matching the saved selection state.
$ws.Range("G4").Select()
